$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "conditioned", 365, 0),
    @(1, "unconditioned", 365, 0),
    @(2, "conditioned", 365, 2),
    @(2, "unconditioned", 365, 0),
    @(3, "conditioned", 365, 0),
    @(3, "unconditioned", 365, 0),
    @(4, "conditioned", 365, 0),
    @(4, "unconditioned", 365, 0),
    @(5, "conditioned", 365, 0),
    @(5, "unconditioned", 365, 0),
    @(6, "conditioned", 365, 2),
    @(6, "unconditioned", 365, 0),
    @(7, "conditioned", 365, 0),
    @(7, "unconditioned", 365, 0),
    @(8, "conditioned", 365, 1),
    @(8, "unconditioned", 365, 0),
    @(9, "conditioned", 365, 0),
    @(9, "unconditioned", 365, 0),
    @(10, "conditioned", 365, 0),
    @(10, "unconditioned", 365, 0),
    @(11, "conditioned", 365, 0),
    @(11, "unconditioned", 365, 1),
    @(12, "conditioned", 365, 0),
    @(12, "unconditioned", 365, 0),
    @(13, "conditioned", 365, 0),
    @(13, "unconditioned", 365, 1),
    @(14, "conditioned", 365, 0),
    @(14, "unconditioned", 365, 0),
    @(15, "unconditioned", 365, 0)
)

$startRow = 466
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

$ws.Range("C176").Select()
$excel.ActiveWindow.Zoom = 192

Write-Output "done"
